$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each target cell is plain text in the source file (inline string, default
# "Normal" style). Some new values are digit-only (e.g. "572.81"), and the
# COM layer auto-detects those as numbers on a bare .Value assignment, which
# would both change the stored type and pick up a float rounding tail (e.g.
# 572.80999999999995) plus a new number-format style. Prefixing with a leading
# apostrophe forces text entry (quote-prefix) like typing in Excel, then resetting
# the Style back to "Normal" clears the quote-prefix style bit so the cell ends up
# identical in type/format to how it started - just with the new text value.

$ws.Range("D2").Value = "'59.065.83"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.13%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.572.28"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.52%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.03%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'572.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +2.30%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'142.72"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.89%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.14%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.594"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.49%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.575.55"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -2.01%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -2.07%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +2.23%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +11.82%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +2.22%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'3.028.89"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.48%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'59.120.03"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.03%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'22.29"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +5.36%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  +2.75%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.578.06"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -1.78%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  +1.17%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'335.53"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.76%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  +1.00%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  +1.37%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.10%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'64.58"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -2.25%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.463"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +7.99%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +0.54%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -2.24%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'7.27"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +0.78%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'0.0₃0780"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +1.64%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  +0.06%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  +0.08%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'160.28"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +3.57%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +0.26%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E35").Value = "'  +0.35%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +2.30%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.874"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -3.95%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -4.58%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +0.44%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +0.96%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'295.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +3.16%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'3.65"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +1.05%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  +0.19%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'131.41"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +11.33%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.0975"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +1.80%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  -1.47%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.0537"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.75%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  +0.06%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  +1.27%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  +2.11%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'1.944.46"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.53%  "
$ws.Range("E51").Style = "Normal"
